# Update the "dSF" column (F) values for specific rows, per repull of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -6
$ws.Range("F6").Value = 4
$ws.Range("F9").Value = 1
$ws.Range("F11").Value = 1
$ws.Range("F13").Value = -9
$ws.Range("F14").Value = -3
$ws.Range("F19").Value = -8
$ws.Range("F21").Value = 4
$ws.Range("F33").Value = 0
